$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "34-21="
$t.Cell(1,2).Range.Text = "41-31="
$t.Cell(1,3).Range.Text = "51-4="
$t.Cell(1,4).Range.Text = "77+13="
$t.Cell(1,5).Range.Text = "78-75="
$t.Cell(2,1).Range.Text = "48+6="
$t.Cell(2,2).Range.Text = "34+13="
$t.Cell(2,3).Range.Text = "44+25="
$t.Cell(2,4).Range.Text = "13+76="
$t.Cell(2,5).Range.Text = "43+7="
$t.Cell(3,1).Range.Text = "34+37="
$t.Cell(3,2).Range.Text = "79-14="
$t.Cell(3,3).Range.Text = "71-2="
$t.Cell(3,4).Range.Text = "50-18="
$t.Cell(3,5).Range.Text = "52-27="
$t.Cell(4,1).Range.Text = "35+24="
$t.Cell(4,2).Range.Text = "66+6="
$t.Cell(4,3).Range.Text = "59-15="
$t.Cell(4,4).Range.Text = "35+13="
$t.Cell(4,5).Range.Text = "81-70="
$t.Cell(5,1).Range.Text = "36+11="
$t.Cell(5,2).Range.Text = "65-49="
$t.Cell(5,3).Range.Text = "28+18="
$t.Cell(5,4).Range.Text = "51+11="
$t.Cell(5,5).Range.Text = "41+27="
$t.Cell(6,1).Range.Text = "6+77="
$t.Cell(6,2).Range.Text = "51-0="
$t.Cell(6,3).Range.Text = "47-40="
$t.Cell(6,4).Range.Text = "0+59="
$t.Cell(6,5).Range.Text = "36+59="
$t.Cell(7,1).Range.Text = "6+56="
$t.Cell(7,2).Range.Text = "90+4="
$t.Cell(7,3).Range.Text = "77-15="
$t.Cell(7,4).Range.Text = "82-15="
$t.Cell(7,5).Range.Text = "87-18="
$t.Cell(8,1).Range.Text = "74-25="
$t.Cell(8,2).Range.Text = "43+40="
$t.Cell(8,3).Range.Text = "20+46="
$t.Cell(8,4).Range.Text = "86-7="
$t.Cell(8,5).Range.Text = "78+19="
$t.Cell(9,1).Range.Text = "14+41="
$t.Cell(9,2).Range.Text = "14+1="
$t.Cell(9,3).Range.Text = "92-38="
$t.Cell(9,4).Range.Text = "44+39="
$t.Cell(9,5).Range.Text = "53+17="
$t.Cell(10,1).Range.Text = "35+8="
$t.Cell(10,2).Range.Text = "82-54="
$t.Cell(10,3).Range.Text = "35+41="
$t.Cell(10,4).Range.Text = "44-6="
$t.Cell(10,5).Range.Text = "87-1="
$t.Cell(11,1).Range.Text = "66+33="
$t.Cell(11,2).Range.Text = "93+4="
$t.Cell(11,3).Range.Text = "90-28="
$t.Cell(11,4).Range.Text = "58+20="
$t.Cell(11,5).Range.Text = "35-24="
$t.Cell(12,1).Range.Text = "50-12="
$t.Cell(12,2).Range.Text = "54+42="
$t.Cell(12,3).Range.Text = "64-46="
$t.Cell(12,4).Range.Text = "50+41="
$t.Cell(12,5).Range.Text = "47-29="
$t.Cell(13,1).Range.Text = "48+7="
$t.Cell(13,2).Range.Text = "8+56="
$t.Cell(13,3).Range.Text = "64-23="
$t.Cell(13,4).Range.Text = "89-71="
$t.Cell(13,5).Range.Text = "68-31="
$t.Cell(14,1).Range.Text = "36-11="
$t.Cell(14,2).Range.Text = "99-74="
$t.Cell(14,3).Range.Text = "63-50="
$t.Cell(14,4).Range.Text = "1+35="
$t.Cell(14,5).Range.Text = "7+0="
$t.Cell(15,1).Range.Text = "79-38="
$t.Cell(15,2).Range.Text = "78+9="
$t.Cell(15,3).Range.Text = "91-79="
$t.Cell(15,4).Range.Text = "92-9="
$t.Cell(15,5).Range.Text = "20+44="
$t.Cell(16,1).Range.Text = "6+26="
$t.Cell(16,2).Range.Text = "73-53="
$t.Cell(16,3).Range.Text = "60+32="
$t.Cell(16,4).Range.Text = "54+31="
$t.Cell(16,5).Range.Text = "98-26="
$t.Cell(17,1).Range.Text = "49+21="
$t.Cell(17,2).Range.Text = "40-24="
$t.Cell(17,3).Range.Text = "9+11="
$t.Cell(17,4).Range.Text = "84+14="
$t.Cell(17,5).Range.Text = "79+5="
$t.Cell(18,1).Range.Text = "11+23="
$t.Cell(18,2).Range.Text = "97-28="
$t.Cell(18,3).Range.Text = "13+44="
$t.Cell(18,4).Range.Text = "96-61="
$t.Cell(18,5).Range.Text = "23+50="
$t.Cell(19,1).Range.Text = "47+32="
$t.Cell(19,2).Range.Text = "0+45="
$t.Cell(19,3).Range.Text = "49+0="
$t.Cell(19,4).Range.Text = "45+54="
$t.Cell(19,5).Range.Text = "51-26="
$t.Cell(20,1).Range.Text = "44-36="
$t.Cell(20,2).Range.Text = "7+62="
$t.Cell(20,3).Range.Text = "23+42="
$t.Cell(20,4).Range.Text = "22+50="
$t.Cell(20,5).Range.Text = "44-10="
